$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (price + 1h volume change columns)
# Column D values are forced to text via a leading apostrophe so Excel
# doesn't reinterpret numeric-looking strings (e.g. '1.003') as numbers,
# which would silently drop significant trailing zeros.

$ws.Range("D2").Value = "'27.300.20"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "'1.708.61"
$ws.Range("E3").Value = "  -0.83%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'223.75"
$ws.Range("E5").Value = "  -2.47%  "
$ws.Range("D6").Value = "'0.5293"
$ws.Range("E6").Value = "  -2.29%  "
$ws.Range("D7").Value = "'1.003"
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "'0.2639"
$ws.Range("E8").Value = "  -4.22%  "
$ws.Range("D9").Value = "'0.06543"
$ws.Range("E9").Value = "  -3.90%  "
$ws.Range("D10").Value = "'20.92"
$ws.Range("E10").Value = "  -2.77%  "
$ws.Range("D11").Value = "'0.07618"
$ws.Range("E11").Value = "  -2.22%  "
$ws.Range("D12").Value = "'4.565"
$ws.Range("E12").Value = "  -3.18%  "
$ws.Range("D13").Value = "'1.714.91"
$ws.Range("E13").Value = "  -0.90%  "
$ws.Range("D14").Value = "'1.946.32"
$ws.Range("E14").Value = "  -0.64%  "
$ws.Range("D15").Value = "'0.5739"
$ws.Range("E15").Value = "  -4.21%  "
$ws.Range("D16").Value = "'0.0₅8192"
$ws.Range("E16").Value = "  -2.54%  "
$ws.Range("D17").Value = "'67.24"
$ws.Range("E17").Value = "  -2.07%  "
$ws.Range("D18").Value = "'27.304.79"
$ws.Range("E18").Value = "  -0.61%  "
$ws.Range("D19").Value = "'215.14"
$ws.Range("E19").Value = "  +2.18%  "
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").Value = "'4.672"
$ws.Range("E21").Value = "  -2.72%  "
$ws.Range("D22").Value = "'10.46"
$ws.Range("E22").Value = "  -4.25%  "
$ws.Range("D23").Value = "'5.968"
$ws.Range("E23").Value = "  -4.06%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Value = "'142.75"
$ws.Range("E25").Value = "  -2.45%  "
$ws.Range("D26").Value = "'1.754"
$ws.Range("E26").Value = "  +8.18%  "
$ws.Range("D27").Value = "'0.1218"
$ws.Range("E27").Value = "  -2.68%  "
$ws.Range("E28").Value = "  -2.23%  "
$ws.Range("D29").Value = "'16.31"
$ws.Range("E29").Value = "  -3.31%  "
$ws.Range("D30").Value = "'0.05377"
$ws.Range("E30").Value = "  -3.92%  "
$ws.Range("D31").Value = "'1.293"
$ws.Range("E31").Value = "  -1.64%  "
$ws.Range("D32").Value = "'3.491"
$ws.Range("E32").Value = "  -4.84%  "
$ws.Range("E33").Value = "  -3.04%  "
$ws.Range("D34").Value = "'1.639"
$ws.Range("E34").Value = "  +0.73%  "
$ws.Range("D35").Value = "'2.872"
$ws.Range("E35").Value = "  +0.75%  "
$ws.Range("D36").Value = "'0.9498"
$ws.Range("E36").Value = "  -2.75%  "
$ws.Range("E37").Value = "  -0.95%  "
$ws.Range("D38").Value = "'0.5873"
$ws.Range("E38").Value = "  +0.28%  "
$ws.Range("D39").Value = "'0.01627"
$ws.Range("E39").Value = "  -1.03%  "
$ws.Range("D40").Value = "'5.870"
$ws.Range("E40").Value = "  +0.55%  "
$ws.Range("D41").Value = "'1.003"
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("D42").Value = "'1.040.45"
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("D43").Value = "'0.8396"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").Value = "'101.03"
$ws.Range("E44").Value = "  -1.38%  "
$ws.Range("D45").Value = "'1.853.28"
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("D47").Value = "'58.05"
$ws.Range("E47").Value = "  -2.63%  "
$ws.Range("D48").Value = "'0.4498"
$ws.Range("E48").Value = "  +1.58%  "
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("E50").Value = "  -1.70%  "
$ws.Range("E51").Value = "  -0.73%  "
